# Generate Report for Handback
#
# Row 7 ("3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2") on both the "zh-cn" and
# "de-de" sheets gets a handback result recorded: the Latest Target File
# (I), Latest Handback File (J) and Latest Handback DateTime (K) columns
# are populated, and an Error Detail (P) is recorded because the handback
# file version was stale.

$wb = $excel.ActiveWorkbook

$mdName        = "3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.md"
$hyperlinkUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a954ad2e8a41affeb1fbe70b7d8e3f56aac13b90/e2e/3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.md"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c3792cb11d91dd9a394c76a333d84e7f0850509/e2e/3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a954ad2e8a41affeb1fbe70b7d8e3f56aac13b90/e2e/3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.md."

function Set-HandbackRow($ws, $handbackXlf, $handbackDateTime) {
    # I7: Latest Target File -> becomes a hyperlink to the handback .md,
    # mirroring the styling already used by A7/I2 (HyperLink cell style).
    $i7 = $ws.Range("I7")
    $ws.Hyperlinks.Add($i7, $hyperlinkUrl, "", "", $mdName)
    $i7.Font.Underline = 2
    $i7.Font.Color = 15570276

    # J7: Latest Handback File
    $ws.Range("J7").Value = $handbackXlf

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # P7: Error Detail
    $ws.Range("P7").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn "3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.d7b6c16da37a597812bca18fa2a4f02510d8c929.zh-cn.xlf" "2016-08-15 16:52:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe "3f64141b-5c1d-48d8-92cf-5e8ec2d68dc2.d7b6c16da37a597812bca18fa2a4f02510d8c929.de-de.xlf" "2016-08-15 16:52:22"
